$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: refresh the daily price records for Cebollín (Terminal
# Hortofrutícola Agro Chillán) with the latest Fecha / Volumen / Precio /
# Unidad de comercialización / Origen / Precio $/Kg / Kg o Unidades values.
# Each element below carries the full record that now belongs to that row.
$data = @(
    @{Row=2; D=44818; J=120; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=3; D=44208; J=85; K=3700; L=4000; M=3824; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1912; Q=2},
    @{Row=4; D=44704; J=100; K=6000; L=6500; M=6250; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=174; Q=36},
    @{Row=5; D=44762; J=60; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región Metropolitana"; P=2667; Q=3},
    @{Row=6; D=44664; J=200; K=8000; L=8500; M=8250; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=229; Q=36},
    @{Row=7; D=44769; J=100; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2667; Q=3},
    @{Row=8; D=44764; J=100; K=8000; L=9000; M=8500; N="`$/docena de atados"; O="Región Metropolitana"; P=2833; Q=3},
    @{Row=9; D=44223; J=80; K=3500; L=3800; M=3688; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1844; Q=2},
    @{Row=10; D=44798; J=200; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=11; D=44760; J=120; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región Metropolitana"; P=2667; Q=3},
    @{Row=12; D=44804; J=120; K=8500; L=9000; M=8750; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2917; Q=3},
    @{Row=13; D=44662; J=200; K=8000; L=8500; M=8250; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=229; Q=36},
    @{Row=14; D=44817; J=120; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=15; D=44701; J=120; K=7000; L=7500; M=7250; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=201; Q=36},
    @{Row=16; D=44810; J=120; K=8000; L=9000; M=8500; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2833; Q=3},
    @{Row=18; D=44790; J=120; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=19; D=44160; J=43; K=3500; L=4000; M=3709; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=103; Q=36},
    @{Row=20; D=44771; J=150; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2667; Q=3},
    @{Row=21; D=44215; J=140; K=3500; L=4000; M=3768; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1884; Q=2},
    @{Row=22; D=44210; J=105; K=3500; L=4000; M=3714; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1857; Q=2},
    @{Row=23; D=44791; J=120; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=24; D=44811; J=100; K=8000; L=9000; M=8500; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2833; Q=3},
    @{Row=25; D=44166; J=70; K=3500; L=4000; M=3679; N="`$/paquete 36 unidades"; O="Región Metropolitana"; P=102; Q=36},
    @{Row=26; D=44225; J=80; K=3400; L=3700; M=3550; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1775; Q=2},
    @{Row=27; D=44161; J=50; K=2800; L=3000; M=2900; N="`$/paquete 2 kilos"; O="Provincia de Diguillín"; P=1450; Q=2},
    @{Row=28; D=44784; J=160; K=8000; L=8500; M=8250; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2750; Q=3},
    @{Row=30; D=44775; J=100; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Provincia de Diguillín"; P=2667; Q=3}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # Fecha
    $ws.Cells.Item($r, 10).Value = $item.J   # Volumen
    $ws.Cells.Item($r, 11).Value = $item.K   # Precio minimo
    $ws.Cells.Item($r, 12).Value = $item.L   # Precio maximo
    $ws.Cells.Item($r, 13).Value = $item.M   # Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $item.N   # Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $item.O   # Origen
    $ws.Cells.Item($r, 16).Value = $item.P   # Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $item.Q   # Kg o Unidades
}
